{"js": "// The document previously held four TODO-style notes (primary key,\n// restrictions, index renaming, and the idioma/genero tables reminder).\n// That last reminder has now been acted on, so all four notes are\n// removed, leaving a single empty paragraph in the body.\nconst body = context.document.body;\n\n// Clear all content from the body; this collapses everything down to\n// one empty trailing paragraph (required - a Word body can never be\n// fully empty).\nbody.clear();\nawait context.sync();\n\n// `clear()` leaves behind a lone empty run inside that paragraph; drop\n// it too so the remaining paragraph is completely empty.\nconst firstPara = body.paragraphs.getFirst();\nfirstPara.getRange().delete();\nawait context.sync();\n", "ps1": "# The document previously held four TODO-style notes (primary key,\n# restrictions, index renaming, and the idioma/genero tables reminder).\n# That last reminder has now been acted on, so all four notes are\n# removed, leaving a single empty paragraph in the body.\n$d = $word.ActiveDocument\n\n# A Word body can never be completely empty, so first append a brand\n# new (attribute-free) blank paragraph after the last existing one -\n# this will become the sole surviving paragraph.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n\n# Remove the four original note paragraphs (always paragraph 1, since\n# each deletion shifts the following ones up).\nfor ($i = 1; $i -le 4; $i++) {\n    $d.Paragraphs(1).Range.Delete()\n}\n\n# The freshly inserted paragraph still carries an empty run; drop it so\n# the remaining paragraph is completely empty.\n$d.Paragraphs(1).Range.Delete()\n"}
